$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 462, shifting existing rows 462-556 down to 464-558
$ws.Range("A462:A463").EntireRow.Insert()

# Fill in the new row 462 (Camote, 1a (cosecha))
$ws.Cells.Item(462, 1).Value = 3
$ws.Cells.Item(462, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(462, 3).Value = "Coquimbo"
$ws.Cells.Item(462, 4).Value = 44711
$ws.Cells.Item(462, 5).Value = 5
$ws.Cells.Item(462, 6).Value = 100112045
$ws.Cells.Item(462, 7).Value = "Zapallo"
$ws.Cells.Item(462, 8).Value = "Camote"
$ws.Cells.Item(462, 9).Value = "1a (cosecha)"
$ws.Cells.Item(462, 10).Value = 280
$ws.Cells.Item(462, 11).Value = 600
$ws.Cells.Item(462, 12).Value = 650
$ws.Cells.Item(462, 13).Value = 629
$ws.Cells.Item(462, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(462, 15).Value = "Provincia de Talca"
$ws.Cells.Item(462, 16).Value = 629
$ws.Cells.Item(462, 17).Value = 1
$ws.Cells.Item(462, 18).Value = "Hortaliza"

# Fill in the new row 463 (Paine, 1a (cosecha))
$ws.Cells.Item(463, 1).Value = 3
$ws.Cells.Item(463, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(463, 3).Value = "Coquimbo"
$ws.Cells.Item(463, 4).Value = 44711
$ws.Cells.Item(463, 5).Value = 5
$ws.Cells.Item(463, 6).Value = 100112045
$ws.Cells.Item(463, 7).Value = "Zapallo"
$ws.Cells.Item(463, 8).Value = "Paine"
$ws.Cells.Item(463, 9).Value = "1a (cosecha)"
$ws.Cells.Item(463, 10).Value = 370
$ws.Cells.Item(463, 11).Value = 350
$ws.Cells.Item(463, 12).Value = 370
$ws.Cells.Item(463, 13).Value = 356
$ws.Cells.Item(463, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(463, 15).Value = "Provincia de Talca"
$ws.Cells.Item(463, 16).Value = 356
$ws.Cells.Item(463, 17).Value = 1
$ws.Cells.Item(463, 18).Value = "Hortaliza"
